$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "57.836.96"
$ws.Range("E2").Value = "  -3.70%  "
$ws.Range("D3").Value = "2.289.83"
$ws.Range("E3").Value = "  -4.47%  "
Set-TextValue $ws "D5" "534.14"
$ws.Range("E5").Value = "  -4.59%  "
Set-TextValue $ws "D6" "131.09"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue $ws "D8" "0.573"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "2.288.08"
$ws.Range("E9").Value = "  -4.50%  "
Set-TextValue $ws "D10" "0.0994"
$ws.Range("E10").Value = "  -5.86%  "
Set-TextValue $ws "D11" "5.43"
$ws.Range("E11").Value = "  -3.84%  "
$ws.Range("E12").Value = "  -0.59%  "
Set-TextValue $ws "D13" "0.331"
$ws.Range("E13").Value = "  -4.19%  "
Set-TextValue $ws "D14" "23.39"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "2.700.07"
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "57.849.49"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").Value = "2.275.98"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("E19").Value = "  -5.41%  "
$ws.Range("E20").Value = "  -6.46%  "
Set-TextValue $ws "D21" "313.36"
$ws.Range("E21").Value = "  -2.98%  "
Set-TextValue $ws "D22" "6.39"
$ws.Range("E22").Value = "  -5.57%  "
Set-TextValue $ws "D23" "0.999"
$ws.Range("E23").Value = "  -0.16%  "
Set-TextValue $ws "D24" "62.76"
$ws.Range("E24").Value = "  -2.20%  "
Set-TextValue $ws "D25" "0.166"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("E28").Value = "  -5.43%  "
Set-TextValue $ws "D29" "169.82"
$ws.Range("E29").Value = "  -0.69%  "
Set-TextValue $ws "D30" "1.70"
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("E31").Value = "  -6.38%  "
$ws.Range("E32").Value = "  -6.11%  "
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("E34").Value = "  -5.78%  "
$ws.Range("E35").Value = "  -0.03%  "
Set-TextValue $ws "D36" "17.68"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("E37").Value = "  -0.01%  "
Set-TextValue $ws "D38" "1.23"
$ws.Range("E38").Value = "  -7.46%  "
$ws.Range("E39").Value = "  -6.72%  "
Set-TextValue $ws "D40" "38.05"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("E41").Value = "  -6.75%  "
Set-TextValue $ws "D42" "141.35"
$ws.Range("E42").Value = "  -4.07%  "
Set-TextValue $ws "D43" "289.36"
$ws.Range("E43").Value = "  -10.56%  "
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("E45").Value = "  -2.49%  "
Set-TextValue $ws "D46" "0.0495"
$ws.Range("E46").Value = "  -3.51%  "
Set-TextValue $ws "D47" "0.556"
$ws.Range("E47").Value = "  -3.36%  "
Set-TextValue $ws "D48" "18.07"
$ws.Range("E48").Value = "  -9.08%  "
$ws.Range("E49").Value = "  -4.51%  "
Set-TextValue $ws "D50" "10.96"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "0.0₆0203"
$ws.Range("E51").Value = "  +83.82%  "
